$d = $word.ActiveDocument

# ------------------------------------------------------------------
# "bitch" -> "Bitch x2"
#   - capitalise the first letter
#   - append " x2" at the end
# Locate the word precisely with Find so we don't depend on fixed
# character offsets.
# ------------------------------------------------------------------
$target = $d.Content
$found = $target.Find.Execute("bitch", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $wordStart = $target.Start
    $wordEnd   = $target.End

    # Capitalise just the first character ("b" -> "B").
    $firstChar = $d.Range($wordStart, $wordStart + 1)
    $firstChar.Text = "B"

    # Append " x2" right after the word ("bitch" -> "bitch x2").
    $tail = $d.Range($wordEnd, $wordEnd)
    $tail.InsertAfter(" x2")
}
